# Cell.xlsx — add a cell whose formula evaluates to an error (#NAME?), entered
# as a (dynamic-)array formula, so the test suite can exercise ingesting
# error-valued cells as 'STRING' Cell's.
#
# Equivalent of a user typing  =l+1  into AN1 and confirming it with
# Ctrl+Shift+Enter (array entry) — "l" isn't a defined name, so Excel raises
# #NAME? and stores it as an array formula over the single cell AN1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$an1 = $ws.Range("AN1")

# Array-enter the formula; resolves to the #NAME? error because "l" is not a
# recognised name/function.
$an1.FormulaArray = "=l+1"

# Typing into AN1 leaves it as the active/selected cell (previously AM1 was
# selected).
$an1.Select() | Out-Null

# The column hosting the new (narrower) content re-sizes from its previous
# best-fit width down to fit "#NAME?".
$ws.Columns.Item(40).ColumnWidth = 7.65
